$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update price values in column D (rows 14-17)
$ws.Range("D14").Value = 98.8
$ws.Range("D15").Value = 142
$ws.Range("D16").Value = 202
$ws.Range("D17").Value = 361
